$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width
$ws.Columns.Item(1).ColumnWidth = 26.59

# Row 2
$ws.Range("A2").Value = "Seguro de Incendio"
$ws.Range("B2").Value = "ADKF123454"
$ws.Range("C2").Value = 15000
$ws.Range("D2").Value = 42.5
$ws.Range("E2").Value = 22
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = "ANGEL"
$ws.Range("H2").Value = 1

# Row 3
$ws.Range("A3").Value = "Seguro de Responsabilidad Civil"
$ws.Range("B3").Value = "AJSM323422"
$ws.Range("C3").Value = 20000
$ws.Range("D3").Value = 20.3
$ws.Range("E3").Value = 18
$ws.Range("F3").Value = 55
$ws.Range("G3").Value = "ANGEL"
$ws.Range("H3").Value = 1

# Row 4
$ws.Range("A4").Value = "Seguro Funerario"
$ws.Range("B4").Value = "Kdma458155"
$ws.Range("C4").Value = 50000
$ws.Range("D4").Value = 35.62
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 99
$ws.Range("G4").Value = "ANGEL"
$ws.Range("H4").Value = 1

# Row 5
$ws.Range("A5").Value = "Seguro Medico Plus"
$ws.Range("B5").Value = "PLSM234333"
$ws.Range("C5").Value = 60000
$ws.Range("D5").Value = 43.6
$ws.Range("E5").Value = 25
$ws.Range("F5").Value = 80
$ws.Range("G5").Value = "ANGEL"
$ws.Range("H5").Value = 1

# Row 6
$ws.Range("A6").Value = "Seguro Automotriz"
$ws.Range("B6").Value = "ATRUD123443"
$ws.Range("C6").Value = 70000
$ws.Range("D6").Value = 60.5
$ws.Range("E6").Value = 18
$ws.Range("F6").Value = 70
$ws.Range("G6").Value = "ANGEL"
$ws.Range("H6").Value = 1

# Update selection to F6
$ws.Range("F6").Select()
